# Scheduled market-price refresh: updates the scraped price/profit columns
# (H:N = currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the leves whose
# underlying item prices moved since the last run. Only specific rows on
# specific sheets are touched; everything else (names, levels, gil, etc.) is
# left exactly as-is.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1111355.6
$ws.Range("I19").Value = 2222400.5
$ws.Range("J19").Value = 310.83334
$ws.Range("K19").Value = 2222400.5
$ws.Range("L19").Value = 310.83334
$ws.Range("M19").Value = -2222225.5
$ws.Range("N19").Value = -660.83334
# Row 62
$ws.Range("H62").Value = 2066.75
$ws.Range("I62").Value = 2055.7778
$ws.Range("J62").Value = 2099.6667
$ws.Range("K62").Value = 2055.7778
$ws.Range("L62").Value = 2099.6667
$ws.Range("M62").Value = -1431.7778
$ws.Range("N62").Value = -3347.6667
# Row 65
$ws.Range("H65").Value = 2066.75
$ws.Range("I65").Value = 2055.7778
$ws.Range("J65").Value = 2099.6667
$ws.Range("K65").Value = 10278.889
$ws.Range("L65").Value = 10498.3335
$ws.Range("M65").Value = -7158.888999999999
$ws.Range("N65").Value = -16738.3335
# Row 111
$ws.Range("H111").Value = 2232.25
$ws.Range("I111").Value = 2143
$ws.Range("J111").Value = 2500
$ws.Range("K111").Value = 6429
$ws.Range("L111").Value = 7500
$ws.Range("M111").Value = -3362
$ws.Range("N111").Value = -13634
# Row 125
$ws.Range("H125").Value = 3900
$ws.Range("J125").Value = 3900
$ws.Range("L125").Value = 35100
$ws.Range("N125").Value = -40020
# Row 129
$ws.Range("H129").Value = 978.0909
$ws.Range("I129").Value = 349.125
$ws.Range("J129").Value = 1179.36
$ws.Range("K129").Value = 1047.375
$ws.Range("L129").Value = 3538.08
$ws.Range("M129").Value = 3952.625
$ws.Range("N129").Value = -13538.08
# Row 135
$ws.Range("H135").Value = 1335
$ws.Range("I135").Value = 1172
$ws.Range("J135").Value = 1742.5
$ws.Range("K135").Value = 10548
$ws.Range("L135").Value = 15682.5
$ws.Range("M135").Value = -8013
$ws.Range("N135").Value = -20752.5
# Row 137
$ws.Range("H137").Value = 3721.6445
$ws.Range("I137").Value = 3111.861
$ws.Range("J137").Value = 6160.778
$ws.Range("K137").Value = 9335.582999999999
$ws.Range("L137").Value = 18482.334
$ws.Range("M137").Value = -6785.582999999999
$ws.Range("N137").Value = -23582.334
# Row 138
$ws.Range("H138").Value = 2273.22
$ws.Range("I138").Value = 1180.5278
$ws.Range("J138").Value = 2887.8594
$ws.Range("K138").Value = 3541.5834
$ws.Range("L138").Value = 8663.5782
$ws.Range("M138").Value = 1598.4166
$ws.Range("N138").Value = -18943.5782

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 56
$ws.Range("H56").Value = 50000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 50000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 50000
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -51484
# Row 61
$ws.Range("H61").Value = 1855.6666
$ws.Range("I61").Value = 1423.2632
$ws.Range("K61").Value = 1423.2632
$ws.Range("M61").Value = -1211.2632
# Row 74
$ws.Range("H74").Value = 1961.12
$ws.Range("I74").Value = 1252.45
$ws.Range("J74").Value = 4795.8
$ws.Range("K74").Value = 1252.45
$ws.Range("L74").Value = 4795.8
$ws.Range("M74").Value = -378.45
$ws.Range("N74").Value = -6543.8
# Row 77
$ws.Range("H77").Value = 1961.12
$ws.Range("I77").Value = 1252.45
$ws.Range("J77").Value = 4795.8
$ws.Range("K77").Value = 6262.25
$ws.Range("L77").Value = 23979
$ws.Range("M77").Value = -1894.25
$ws.Range("N77").Value = -32715
# Row 132
$ws.Range("H132").Value = 2917.7917
$ws.Range("I132").Value = 1459.25
$ws.Range("K132").Value = 4377.75
$ws.Range("M132").Value = -1847.75
# Row 136
$ws.Range("H136").Value = 1855.6666
$ws.Range("I136").Value = 1423.2632
$ws.Range("K136").Value = 4269.7896
$ws.Range("M136").Value = -1719.7896
# Row 137
$ws.Range("H137").Value = 53780
$ws.Range("J137").Value = 53780
$ws.Range("L137").Value = 53780
$ws.Range("N137").Value = -63980

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3519.6667
$ws.Range("I31").Value = 1353.9166
$ws.Range("J31").Value = 6407.3335
$ws.Range("K31").Value = 1353.9166
$ws.Range("L31").Value = 6407.3335
$ws.Range("M31").Value = -1058.9166
$ws.Range("N31").Value = -6997.3335
# Row 34
$ws.Range("H34").Value = 3519.6667
$ws.Range("I34").Value = 1353.9166
$ws.Range("J34").Value = 6407.3335
$ws.Range("K34").Value = 1353.9166
$ws.Range("L34").Value = 6407.3335
$ws.Range("M34").Value = -1151.9166
$ws.Range("N34").Value = -6811.3335
# Row 58
$ws.Range("H58").Value = 2237.4465
$ws.Range("I58").Value = 1905.7959
$ws.Range("J58").Value = 4559
$ws.Range("K58").Value = 1905.7959
$ws.Range("L58").Value = 4559
$ws.Range("M58").Value = -1702.7959
$ws.Range("N58").Value = -4965
# Row 134
$ws.Range("H134").Value = 6215.9614
$ws.Range("I134").Value = 6710.85
$ws.Range("J134").Value = 4566.3335
$ws.Range("K134").Value = 20132.55
$ws.Range("L134").Value = 13699.0005
$ws.Range("M134").Value = -17597.55
$ws.Range("N134").Value = -18769.0005
# Row 136
$ws.Range("H136").Value = 2237.4465
$ws.Range("I136").Value = 1905.7959
$ws.Range("J136").Value = 4559
$ws.Range("K136").Value = 5717.3877
$ws.Range("L136").Value = 13677
$ws.Range("M136").Value = -3167.3877
$ws.Range("N136").Value = -18777

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 57
$ws.Range("H57").Value = 1399.1666
$ws.Range("I57").Value = 1399.5
$ws.Range("J57").Value = 1399
$ws.Range("K57").Value = 4198.5
$ws.Range("L57").Value = 4197
$ws.Range("M57").Value = -3639.5
$ws.Range("N57").Value = -5315
# Row 117
$ws.Range("H117").Value = 1421.4286
$ws.Range("I117").Value = 1158.3334
$ws.Range("K117").Value = 3475.0002
$ws.Range("M117").Value = -33.00019999999995
# Row 129
$ws.Range("H129").Value = 2397.1667
$ws.Range("I129").Value = 3407
$ws.Range("J129").Value = 2008.7693
$ws.Range("K129").Value = 10221
$ws.Range("L129").Value = 6026.3079
$ws.Range("M129").Value = -5221
$ws.Range("N129").Value = -16026.3079
# Row 131
$ws.Range("H131").Value = 10639306
$ws.Range("I131").Value = 22728216
$ws.Range("J131").Value = 1064.12
$ws.Range("K131").Value = 68184648
$ws.Range("L131").Value = 3192.36
$ws.Range("M131").Value = -68179608
$ws.Range("N131").Value = -13272.36
# Row 132
$ws.Range("H132").Value = 3137.9524
$ws.Range("I132").Value = 1299.2
$ws.Range("J132").Value = 4809.5454
$ws.Range("K132").Value = 11692.8
$ws.Range("L132").Value = 43285.9086
$ws.Range("M132").Value = -9162.800000000001
$ws.Range("N132").Value = -48345.9086

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3902.3186
$ws.Range("I126").Value = 2816.7083
$ws.Range("J126").Value = 5114.1626
$ws.Range("K126").Value = 8450.124899999999
$ws.Range("L126").Value = 15342.4878
$ws.Range("M126").Value = -5980.124899999999
$ws.Range("N126").Value = -20282.4878
# Row 132
$ws.Range("H132").Value = 2200.1667
$ws.Range("I132").Value = 1075.5238
$ws.Range("K132").Value = 3226.5714
$ws.Range("M132").Value = -696.5713999999998

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4923
$ws.Range("I7").Value = 3800.6667
$ws.Range("K7").Value = 3800.6667
$ws.Range("M7").Value = -3688.6667
# Row 13
$ws.Range("H13").Value = 11151.5
$ws.Range("I13").Value = 206
$ws.Range("K13").Value = 206
$ws.Range("M13").Value = -66
# Row 126
$ws.Range("H126").Value = 4923
$ws.Range("I126").Value = 3800.6667
$ws.Range("K126").Value = 11402.0001
$ws.Range("M126").Value = -8932.000100000001

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 55996.5
$ws.Range("J2").Value = 55996.5
$ws.Range("L2").Value = 55996.5
$ws.Range("N2").Value = -56220.5
# Row 37
$ws.Range("H37").Value = 16030
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 16030
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 16030
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -16436

